$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "x size" / "y size" columns (N, O) to the settings table.
$ws.Range("N1").Value = "x size"
$ws.Range("O1").Value = "y size"

$ws.Range("N2").Value = "float"
$ws.Range("O2").Value = "float"

$ws.Range("N3").Value = 100
$ws.Range("O3").Value = 100

$ws.Range("N4").Value = 10000
$ws.Range("O4").Value = 10000

# Match the selection state recorded in the saved workbook.
$ws.Range("O2").Select()
